# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 834
$wsExhibit.Range("F5").Value = 1032
$wsExhibit.Range("F6").Value = 2402
$wsExhibit.Range("F7").Value = 203

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 834
$wsAll.Range("F7").Value = 1032
$wsAll.Range("F8").Value = 2403
$wsAll.Range("F10").Value = 203
